$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 504). The whole column was bumped from 45172 (2023-09-03)
# to 45175 (2023-09-06). Update the values in place; the existing cell
# style/format (s="1", date format) is left untouched since we are only
# writing the underlying numeric value.
$ws.Range("C2:C504").Value = 45175
